$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): a new "Num" column is inserted after the pk column,
#     and a new "Name" column is inserted before GachaProbNum. ---
$ws.Range("A1").Value = "Num"
$ws.Range("B1").Value = "Num"
$ws.Range("C1").Value = "Tag"
$ws.Range("D1").Value = "Order"
$ws.Range("E1").Value = "DisplayOrder"
$ws.Range("F1").Value = "Name"
$ws.Range("G1").Value = "GachaProbNum"
$ws.Range("H1").Value = "PickupCookieNum"

# --- Type row (row 2) ---
$ws.Range("A2").Value = "int:pk"
$ws.Range("B2").Value = "int"
$ws.Range("C2").Value = "string"
$ws.Range("D2").Value = "int"
$ws.Range("E2").Value = "int"
$ws.Range("F2").Value = "string"
$ws.Range("G2").Value = "int"
$ws.Range("H2").Value = "int"

# --- Annotation row (row 3) stays put (not affected by the new columns) ---
$ws.Range("A3").Value = "#dd"
$ws.Range("B3").Value = "#"

# --- Data row 4 (existing COOKIE_NORMAL row, now with a Num pk + Name) ---
$ws.Range("A4").Value = 10250219
$ws.Range("B4").Value = 1001001
$ws.Range("C4").Value = "COOKIE_NORMAL"
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = "일반 뽑기"
$ws.Range("G4").Value = 1001
$ws.Range("H4").Value = 0

# --- New data row 5 (COOKIE_PICKUP_1) ---
$ws.Range("B5").Value = 1002001
$ws.Range("C5").Value = "COOKIE_PICKUP_1"
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = "픽업 뽑기"
$ws.Range("G5").Value = 2001
$ws.Range("H5").Value = 6001

# --- Column widths: B through F share the 18.140625-wide formatting that used
#     to cover B:D (column E, which used to be the narrower GachaProbNum
#     column, now holds DisplayOrder and widens to match); G inherits the
#     17.42578125 width that E used to have (GachaProbNum lives there now);
#     H is a brand new, wider column for the PickupCookieNum values. ---
$narrowWidth = $ws.Columns.Item(5).ColumnWidth
$ws.Columns.Item(2).ColumnWidth = $ws.Columns.Item(3).ColumnWidth
$ws.Columns.Item(5).ColumnWidth = $ws.Columns.Item(3).ColumnWidth
$ws.Columns.Item(6).ColumnWidth = $ws.Columns.Item(3).ColumnWidth
$ws.Columns.Item(7).ColumnWidth = $narrowWidth
$ws.Columns.Item(8).ColumnWidth = 19.0

$ws.Range("H4").Select() | Out-Null
